# Update the "Förändrad" (Changed) date column (C) from serial date 46060
# to 46061 for all data rows (rows 2 through 14) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
